$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize header cells C5/D5 to match A5/B5 formatting (same bold style)
$ws.Range("C5:D5").Font.Bold = $true

# Add new flight-data columns: Meal, Credit Card Number
$ws.Range("E5").Value = "Meal"
$ws.Range("F5").Value = "Credit Card Number"
$ws.Range("E5:F5").Font.Bold = $true

# Add new data row values
$ws.Range("E6").Value = "Hindu"
$ws.Range("F6").Value = 122200000000

# Widen column F to fit new header text (closest achievable to the
# author's 18.66 given this engine's internal width quantization)
$ws.Columns.Item(6).ColumnWidth = 17.8

# Move active selection to G7 (as left by the author after editing)
$ws.Range("G7").Select() | Out-Null
